$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 35.995988
$ws.Range("H2").Value = 107.987964
$ws.Range("I2").Value = 0.5613901502831141
$ws.Range("J2").Value = 0.561390150283114
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 34.07074633333333
$ws.Range("N2").Value = 102.212239
$ws.Range("O2").Value = 0.5171464495142372
$ws.Range("P2").Value = 0.5171464495142373
$ws.Range("Q2").Value = 1226.41017616571
$ws.Range("R2").Value = 11037.6915854914
$ws.Range("S2").Value = 0.2903209230111765
$ws.Range("T2").Value = 0.2903209230111765

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 35.995988
$ws.Range("H3").Value = 107.987964
$ws.Range("I3").Value = 0.5613901502831141
$ws.Range("J3").Value = 0.561390150283114
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.685497
$ws.Range("N3").Value = 83.056491
$ws.Range("O3").Value = 0.420227262899125
$ws.Range("P3").Value = 0.4202272628991251
$ws.Range("Q3").Value = 996.566817786036
$ws.Range("R3").Value = 8969.101360074324
$ws.Range("S3").Value = 0.2359114462720015
$ws.Range("T3").Value = 0.2359114462720015

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 35.995988
$ws.Range("H4").Value = 107.987964
$ws.Range("I4").Value = 0.5613901502831141
$ws.Range("J4").Value = 0.561390150283114
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.125957666666666
$ws.Range("N4").Value = 12.377873
$ws.Range("O4").Value = 0.06262628758663766
$ws.Range("P4").Value = 0.06262628758663766
$ws.Range("Q4").Value = 148.5179226578414
$ws.Range("R4").Value = 1336.661303920572
$ws.Range("S4").Value = 0.03515778099993604
$ws.Range("T4").Value = 0.03515778099993603

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 20.666474
$ws.Range("H5").Value = 61.999422
$ws.Range("I5").Value = 0.3223124461726698
$ws.Range("J5").Value = 0.3223124461726698
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 34.07074633333333
$ws.Range("N5").Value = 102.212239
$ws.Range("O5").Value = 0.5171464495142372
$ws.Range("P5").Value = 0.5171464495142373
$ws.Range("Q5").Value = 704.1221932584284
$ws.Range("R5").Value = 6337.099739325857
$ws.Range("S5").Value = 0.1666827371724449
$ws.Range("T5").Value = 0.1666827371724449

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 20.666474
$ws.Range("H6").Value = 61.999422
$ws.Range("I6").Value = 0.3223124461726698
$ws.Range("J6").Value = 0.3223124461726698
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 27.685497
$ws.Range("N6").Value = 83.056491
$ws.Range("O6").Value = 0.420227262899125
$ws.Range("P6").Value = 0.4202272628991251
$ws.Range("Q6").Value = 572.1616039275779
$ws.Range("R6").Value = 5149.454435348202
$ws.Range("S6").Value = 0.1354444770534626
$ws.Range("T6").Value = 0.1354444770534626

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 20.666474
$ws.Range("H7").Value = 61.999422
$ws.Range("I7").Value = 0.3223124461726698
$ws.Range("J7").Value = 0.3223124461726698
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.125957666666666
$ws.Range("N7").Value = 12.377873
$ws.Range("O7").Value = 0.06262628758663766
$ws.Range("P7").Value = 0.06262628758663766
$ws.Range("Q7").Value = 85.26899684326732
$ws.Range("R7").Value = 767.4209715894059
$ws.Range("S7").Value = 0.02018523194676229
$ws.Range("T7").Value = 0.02018523194676229

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 7.456917333333333
$ws.Range("H8").Value = 22.370752
$ws.Range("I8").Value = 0.116297403544216
$ws.Range("J8").Value = 0.116297403544216
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 34.07074633333333
$ws.Range("N8").Value = 102.212239
$ws.Range("O8").Value = 0.5171464495142372
$ws.Range("P8").Value = 0.5171464495142373
$ws.Range("Q8").Value = 254.0627388926364
$ws.Range("R8").Value = 2286.564650033728
$ws.Range("S8").Value = 0.06014278933061578
$ws.Range("T8").Value = 0.06014278933061579

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 7.456917333333333
$ws.Range("H9").Value = 22.370752
$ws.Range("I9").Value = 0.116297403544216
$ws.Range("J9").Value = 0.116297403544216
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 27.685497
$ws.Range("N9").Value = 83.056491
$ws.Range("O9").Value = 0.420227262899125
$ws.Range("P9").Value = 0.4202272628991251
$ws.Range("Q9").Value = 206.448462461248
$ws.Range("R9").Value = 1858.036162151232
$ws.Range("S9").Value = 0.0488713395736609
$ws.Range("T9").Value = 0.04887133957366091

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 7.456917333333333
$ws.Range("H10").Value = 22.370752
$ws.Range("I10").Value = 0.116297403544216
$ws.Range("J10").Value = 0.116297403544216
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.125957666666666
$ws.Range("N10").Value = 12.377873
$ws.Range("O10").Value = 0.06262628758663766
$ws.Range("P10").Value = 0.06262628758663766
$ws.Range("Q10").Value = 30.76692524116622
$ws.Range("R10").Value = 276.902327170496
$ws.Range("S10").Value = 0.007283274639939327
$ws.Range("T10").Value = 0.007283274639939327
